$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill F57:F62 with "T" to match the pattern already used in column F for rows 2-56
$ws.Range("F57:F62").Value = "T"

# Update the view to reflect the scroll position / selection used when the
# edit was made (selection F56:F62, active cell F56)
$excel.ActiveWindow.ScrollRow = 37
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F56:F62").Select()

$wb.Save()
